$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds weekly Tuna price records for
# "Mercado Mayorista Lo Valledor de Santiago" ending at row 220
# (dimension A1:T220). A new week of data (2021-10-05, serial 44474)
# needs to be inserted right after row 207, pushing the existing
# rows 208:220 down to 211:223 (dimension becomes A1:T223).

# Insert three blank rows at 208:210 - this shifts rows 208:220
# (and their formatting, incl. the date-formatted column D) down to 211:223.
$ws.Range("A208:T210").EntireRow.Insert()

# --- Row 208: Tuna, calidad "Especial", new week 2021-10-05 ---
$ws.Range("A208").Value = 6
$ws.Range("B208").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C208").Value = "Metropolitana"
$ws.Range("D208").Value = 44474
$ws.Range("E208").Value = 13
$ws.Range("F208").Value = "Fruta"
$ws.Range("G208").Value = 100107
$ws.Range("H208").Value = "Otros"
$ws.Range("I208").Value = 100107011
$ws.Range("J208").Value = "Tuna"
$ws.Range("K208").Value = "Sin especificar"
$ws.Range("L208").Value = "Especial"
$ws.Range("M208").Value = 210
$ws.Range("N208").Value = 32000
$ws.Range("O208").Value = 32000
$ws.Range("P208").Value = 32000
$ws.Range("Q208").Value = "$/caja 18 kilos"
$ws.Range("R208").Value = "Provincia de Melipilla"
$ws.Range("S208").Value = 1778
$ws.Range("T208").Value = 18

# --- Row 209: Tuna, calidad "Primera", new week 2021-10-05 ---
$ws.Range("A209").Value = 6
$ws.Range("B209").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C209").Value = "Metropolitana"
$ws.Range("D209").Value = 44474
$ws.Range("E209").Value = 13
$ws.Range("F209").Value = "Fruta"
$ws.Range("G209").Value = 100107
$ws.Range("H209").Value = "Otros"
$ws.Range("I209").Value = 100107011
$ws.Range("J209").Value = "Tuna"
$ws.Range("K209").Value = "Sin especificar"
$ws.Range("L209").Value = "Primera"
$ws.Range("M209").Value = 30
$ws.Range("N209").Value = 25000
$ws.Range("O209").Value = 25000
$ws.Range("P209").Value = 25000
$ws.Range("Q209").Value = "$/caja 18 kilos"
$ws.Range("R209").Value = "Provincia de Melipilla"
$ws.Range("S209").Value = 1389
$ws.Range("T209").Value = 18

# --- Row 210: Tuna, calidad "Segunda", new week 2021-10-05 ---
$ws.Range("A210").Value = 6
$ws.Range("B210").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C210").Value = "Metropolitana"
$ws.Range("D210").Value = 44474
$ws.Range("E210").Value = 13
$ws.Range("F210").Value = "Fruta"
$ws.Range("G210").Value = 100107
$ws.Range("H210").Value = "Otros"
$ws.Range("I210").Value = 100107011
$ws.Range("J210").Value = "Tuna"
$ws.Range("K210").Value = "Sin especificar"
$ws.Range("L210").Value = "Segunda"
$ws.Range("M210").Value = 10
$ws.Range("N210").Value = 15000
$ws.Range("O210").Value = 15000
$ws.Range("P210").Value = 15000
$ws.Range("Q210").Value = "$/caja 18 kilos"
$ws.Range("R210").Value = "Provincia de Melipilla"
$ws.Range("S210").Value = 833
$ws.Range("T210").Value = 18
